$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col4a1"
$ws.Range("C2").Value = "Itgb8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 228.2871476666667
$ws.Range("H2").Value = 684.861443
$ws.Range("I2").Value = 0.6105798777018374
$ws.Range("J2").Value = 0.6105798777018375
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.009519999999999999
$ws.Range("N2").Value = 0.02856
$ws.Range("O2").Value = 0.0009583584527718872
$ws.Range("P2").Value = 0.0009583584527718872
$ws.Range("Q2").Value = 2.173293645786667
$ws.Range("R2").Value = 19.55964281208
$ws.Range("S2").Value = 0.000585154386887981
$ws.Range("T2").Value = 0.0005851543868879811

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col4a1"
$ws.Range("C3").Value = "Itgb8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 228.2871476666667
$ws.Range("H3").Value = 684.861443
$ws.Range("I3").Value = 0.6105798777018374
$ws.Range("J3").Value = 0.6105798777018375
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.079146666666667
$ws.Range("N3").Value = 3.23744
$ws.Range("O3").Value = 0.1086354338004839
$ws.Range("P3").Value = 0.1086354338004838
$ws.Range("Q3").Value = 246.3553144473244
$ws.Range("R3").Value = 2217.19783002592
$ws.Range("S3").Value = 0.0663306098839855
$ws.Range("T3").Value = 0.0663306098839855

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col4a1"
$ws.Range("C4").Value = "Itgb8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 228.2871476666667
$ws.Range("H4").Value = 684.861443
$ws.Range("I4").Value = 0.6105798777018374
$ws.Range("J4").Value = 0.6105798777018375
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.844986
$ws.Range("N4").Value = 26.534958
$ws.Range("O4").Value = 0.8904062077467443
$ws.Range("P4").Value = 0.8904062077467442
$ws.Range("Q4").Value = 2019.1966250916
$ws.Range("R4").Value = 18172.7696258244
$ws.Range("S4").Value = 0.543664113430964
$ws.Range("T4").Value = 0.543664113430964

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col4a1"
$ws.Range("C5").Value = "Itgb8"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 56.66021733333333
$ws.Range("H5").Value = 169.980652
$ws.Range("I5").Value = 0.1515441798784964
$ws.Range("J5").Value = 0.1515441798784964
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009519999999999999
$ws.Range("N5").Value = 0.02856
$ws.Range("O5").Value = 0.0009583584527718872
$ws.Range("P5").Value = 0.0009583584527718872
$ws.Range("Q5").Value = 0.5394052690133332
$ws.Range("R5").Value = 4.854647421119999
$ws.Range("S5").Value = 0.0001452336457549403
$ws.Range("T5").Value = 0.0001452336457549404

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col4a1"
$ws.Range("C6").Value = "Itgb8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 56.66021733333333
$ws.Range("H6").Value = 169.980652
$ws.Range("I6").Value = 0.1515441798784964
$ws.Range("J6").Value = 0.1515441798784964
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.079146666666667
$ws.Range("N6").Value = 3.23744
$ws.Range("O6").Value = 0.1086354338004839
$ws.Range("P6").Value = 0.1086354338004838
$ws.Range("Q6").Value = 61.14468466787555
$ws.Range("R6").Value = 550.30216201088
$ws.Range("S6").Value = 0.01646306772103901
$ws.Range("T6").Value = 0.01646306772103901

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col4a1"
$ws.Range("C7").Value = "Itgb8"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 56.66021733333333
$ws.Range("H7").Value = 169.980652
$ws.Range("I7").Value = 0.1515441798784964
$ws.Range("J7").Value = 0.1515441798784964
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.844986
$ws.Range("N7").Value = 26.534958
$ws.Range("O7").Value = 0.8904062077467443
$ws.Range("P7").Value = 0.8904062077467442
$ws.Range("Q7").Value = 501.1588290702907
$ws.Range("R7").Value = 4510.429461632616
$ws.Range("S7").Value = 0.1349358785117024
$ws.Range("T7").Value = 0.1349358785117024

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Col4a1"
$ws.Range("C8").Value = "Itgb8"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6468253333333333
$ws.Range("H8").Value = 1.940476
$ws.Range("I8").Value = 0.001730007742257072
$ws.Range("J8").Value = 0.001730007742257073
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009519999999999999
$ws.Range("N8").Value = 0.02856
$ws.Range("O8").Value = 0.0009583584527718872
$ws.Range("P8").Value = 0.0009583584527718872
$ws.Range("Q8").Value = 0.006157777173333332
$ws.Range("R8").Value = 0.05541999455999999
$ws.Range("S8").Value = 0.000001657967543152874
$ws.Range("T8").Value = 0.000001657967543152874

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Col4a1"
$ws.Range("C9").Value = "Itgb8"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6468253333333333
$ws.Range("H9").Value = 1.940476
$ws.Range("I9").Value = 0.001730007742257072
$ws.Range("J9").Value = 0.001730007742257073
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.079146666666667
$ws.Range("N9").Value = 3.23744
$ws.Range("O9").Value = 0.1086354338004839
$ws.Range("P9").Value = 0.1086354338004838
$ws.Range("Q9").Value = 0.6980194023822222
$ws.Range("R9").Value = 6.28217462144
$ws.Range("S9").Value = 0.0001879401415582927
$ws.Range("T9").Value = 0.0001879401415582927

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Col4a1"
$ws.Range("C10").Value = "Itgb8"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6468253333333333
$ws.Range("H10").Value = 1.940476
$ws.Range("I10").Value = 0.001730007742257072
$ws.Range("J10").Value = 0.001730007742257073
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.844986
$ws.Range("N10").Value = 26.534958
$ws.Range("O10").Value = 0.8904062077467443
$ws.Range("P10").Value = 0.8904062077467442
$ws.Range("Q10").Value = 5.721161017778666
$ws.Range("R10").Value = 51.490449160008
$ws.Range("S10").Value = 0.001540409633155627
$ws.Range("T10").Value = 0.001540409633155627

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col4a1"
$ws.Range("C11").Value = "Itgb8"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7467493333333334
$ws.Range("H11").Value = 2.240248
$ws.Range("I11").Value = 0.00199726581754988
$ws.Range("J11").Value = 0.001997265817549881
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.009519999999999999
$ws.Range("N11").Value = 0.02856
$ws.Range("O11").Value = 0.0009583584527718872
$ws.Range("P11").Value = 0.0009583584527718872
$ws.Range("Q11").Value = 0.007109053653333333
$ws.Range("R11").Value = 0.06398148288000001
$ws.Range("S11").Value = 0.000001914096578681282
$ws.Range("T11").Value = 0.000001914096578681283

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col4a1"
$ws.Range("C12").Value = "Itgb8"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7467493333333334
$ws.Range("H12").Value = 2.240248
$ws.Range("I12").Value = 0.00199726581754988
$ws.Range("J12").Value = 0.001997265817549881
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.079146666666667
$ws.Range("N12").Value = 3.23744
$ws.Range("O12").Value = 0.1086354338004839
$ws.Range("P12").Value = 0.1086354338004838
$ws.Range("Q12").Value = 0.8058520539022223
$ws.Range("R12").Value = 7.252668485120002
$ws.Range("S12").Value = 0.0002169738385044093
$ws.Range("T12").Value = 0.0002169738385044094

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col4a1"
$ws.Range("C13").Value = "Itgb8"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7467493333333334
$ws.Range("H13").Value = 2.240248
$ws.Range("I13").Value = 0.00199726581754988
$ws.Range("J13").Value = 0.001997265817549881
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 8.844986
$ws.Range("N13").Value = 26.534958
$ws.Range("O13").Value = 0.8904062077467443
$ws.Range("P13").Value = 0.8904062077467442
$ws.Range("Q13").Value = 6.604987398842668
$ws.Range("R13").Value = 59.44488658958402
$ws.Range("S13").Value = 0.00177837788246679
$ws.Range("T13").Value = 0.00177837788246679

# Row 14
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Col4a1"
$ws.Range("C14").Value = "Itgb8"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8022086666666667
$ws.Range("H14").Value = 2.406626
$ws.Range("I14").Value = 0.002145598096919091
$ws.Range("J14").Value = 0.002145598096919091
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.009519999999999999
$ws.Range("N14").Value = 0.02856
$ws.Range("O14").Value = 0.0009583584527718872
$ws.Range("P14").Value = 0.0009583584527718872
$ws.Range("Q14").Value = 0.007637026506666666
$ws.Range("R14").Value = 0.06873323856
$ws.Range("S14").Value = 0.000002056252072433685
$ws.Range("T14").Value = 0.000002056252072433686

# Row 15
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Col4a1"
$ws.Range("C15").Value = "Itgb8"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8022086666666667
$ws.Range("H15").Value = 2.406626
$ws.Range("I15").Value = 0.002145598096919091
$ws.Range("J15").Value = 0.002145598096919091
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.079146666666667
$ws.Range("N15").Value = 3.23744
$ws.Range("O15").Value = 0.1086354338004839
$ws.Range("P15").Value = 0.1086354338004838
$ws.Range("Q15").Value = 0.8657008086044445
$ws.Range("R15").Value = 7.791307277440001
$ws.Range("S15").Value = 0.000233087980020298
$ws.Range("T15").Value = 0.0002330879800202981

# Row 16
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Col4a1"
$ws.Range("C16").Value = "Itgb8"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8022086666666667
$ws.Range("H16").Value = 2.406626
$ws.Range("I16").Value = 0.002145598096919091
$ws.Range("J16").Value = 0.002145598096919091
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 8.844986
$ws.Range("N16").Value = 26.534958
$ws.Range("O16").Value = 0.8904062077467443
$ws.Range("P16").Value = 0.8904062077467442
$ws.Range("Q16").Value = 7.095524425745334
$ws.Range("R16").Value = 63.85971983170801
$ws.Range("S16").Value = 0.001910453864826359
$ws.Range("T16").Value = 0.001910453864826359

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col4a1"
$ws.Range("C17").Value = "Itgb8"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 86.74265433333333
$ws.Range("H17").Value = 260.227963
$ws.Range("I17").Value = 0.2320030707629401
$ws.Range("J17").Value = 0.2320030707629402
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.009519999999999999
$ws.Range("N17").Value = 0.02856
$ws.Range("O17").Value = 0.0009583584527718872
$ws.Range("P17").Value = 0.0009583584527718872
$ws.Range("Q17").Value = 0.8257900692533332
$ws.Range("R17").Value = 7.432110623279999
$ws.Range("S17").Value = 0.0002223421039346979
$ws.Range("T17").Value = 0.000222342103934698

# Row 18
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Col4a1"
$ws.Range("C18").Value = "Itgb8"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 86.74265433333333
$ws.Range("H18").Value = 260.227963
$ws.Range("I18").Value = 0.2320030707629401
$ws.Range("J18").Value = 0.2320030707629402
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 1.079146666666667
$ws.Range("N18").Value = 3.23744
$ws.Range("O18").Value = 0.1086354338004839
$ws.Range("P18").Value = 0.1086354338004838
$ws.Range("Q18").Value = 93.60804628163557
$ws.Range("R18").Value = 842.4724165347201
$ws.Range("S18").Value = 0.02520375423537635
$ws.Range("T18").Value = 0.02520375423537636

# Row 19
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Col4a1"
$ws.Range("C19").Value = "Itgb8"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 86.74265433333333
$ws.Range("H19").Value = 260.227963
$ws.Range("I19").Value = 0.2320030707629401
$ws.Range("J19").Value = 0.2320030707629402
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 8.844986
$ws.Range("N19").Value = 26.534958
$ws.Range("O19").Value = 0.8904062077467443
$ws.Range("P19").Value = 0.8904062077467442
$ws.Range("Q19").Value = 767.2375631811727
$ws.Range("R19").Value = 6905.138068630555
$ws.Range("S19").Value = 0.2065769744236291
$ws.Range("T19").Value = 0.2065769744236291
